# The "おめでとう" (mabruuk) post occupying row 563 was removed from the
# posts log. Deleting the entire row shifts every subsequent row (564-687)
# up by one, which also contracts the sheet's used range from C687 to C686 -
# exactly matching the target diff (no cell content elsewhere actually
# changes; only its row number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(563).Delete()
